$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 315.91666
$ws.Range("I33").Value = 359.1
$ws.Range("K33").Value = 359.1
$ws.Range("M33").Value = -130.1
$ws.Range("H47").Value = 46100.5
$ws.Range("I47").Value = 3044.6667
$ws.Range("J47").Value = 89156.336
$ws.Range("K47").Value = 3044.6667
$ws.Range("L47").Value = 89156.336
$ws.Range("M47").Value = -2072.6667
$ws.Range("N47").Value = -91100.336
$ws.Range("H125").Value = 3009.7778
$ws.Range("I125").Value = 3082.3333
$ws.Range("J125").Value = 2864.6667
$ws.Range("K125").Value = 27740.9997
$ws.Range("L125").Value = 25782.0003
$ws.Range("M125").Value = -25280.9997
$ws.Range("N125").Value = -30702.0003
$ws.Range("H133").Value = 74543
$ws.Range("J133").Value = 74543
$ws.Range("L133").Value = 74543
$ws.Range("N133").Value = -84663
$ws.Range("H134").Value = 44993.4
$ws.Range("J134").Value = 44993.4
$ws.Range("L134").Value = 44993.4
$ws.Range("N134").Value = -55133.4
$ws.Range("H135").Value = 1215.8182
$ws.Range("I135").Value = 897.5263
$ws.Range("J135").Value = 3231.6667
$ws.Range("K135").Value = 8077.736699999999
$ws.Range("L135").Value = 29085.0003
$ws.Range("M135").Value = -5542.736699999999
$ws.Range("N135").Value = -34155.0003
$ws.Range("H136").Value = 86776.60000000001
$ws.Range("J136").Value = 86776.60000000001
$ws.Range("L136").Value = 86776.60000000001
$ws.Range("N136").Value = -96976.60000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 821
$ws.Range("I35").Value = 821
$ws.Range("K35").Value = 821
$ws.Range("M35").Value = -415
$ws.Range("H61").Value = 38769.297
$ws.Range("I61").Value = 1720.88
$ws.Range("K61").Value = 1720.88
$ws.Range("M61").Value = -1508.88
$ws.Range("H97").Value = 1076.3077
$ws.Range("I97").Value = 837.3333
$ws.Range("J97").Value = 1614
$ws.Range("K97").Value = 837.3333
$ws.Range("L97").Value = 1614
$ws.Range("M97").Value = -341.3333
$ws.Range("N97").Value = -2606
$ws.Range("H102").Value = 63714.168
$ws.Range("I102").Value = 84720.586
$ws.Range("K102").Value = 84720.586
$ws.Range("M102").Value = -83098.586
$ws.Range("H132").Value = 1494.6909
$ws.Range("I132").Value = 1398.25
$ws.Range("K132").Value = 4194.75
$ws.Range("M132").Value = -1664.75
$ws.Range("H136").Value = 38769.297
$ws.Range("I136").Value = 1720.88
$ws.Range("K136").Value = 5162.64
$ws.Range("M136").Value = -2612.64

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 99990
$ws.Range("J52").Value = 99990
$ws.Range("L52").Value = 99990
$ws.Range("N52").Value = -100516
$ws.Range("H99").Value = 1699082.1
$ws.Range("I99").Value = 72825.78999999999
$ws.Range("K99").Value = 72825.78999999999
$ws.Range("M99").Value = -71327.78999999999
$ws.Range("H109").Value = 77996.664
$ws.Range("J109").Value = 77996.664
$ws.Range("L109").Value = 77996.664
$ws.Range("N109").Value = -80770.664
$ws.Range("H121").Value = 99990
$ws.Range("J121").Value = 99990
$ws.Range("L121").Value = 99990
$ws.Range("N121").Value = -103484
$ws.Range("H132").Value = 36843.125
$ws.Range("J132").Value = 36843.125
$ws.Range("L132").Value = 36843.125
$ws.Range("N132").Value = -46963.125
$ws.Range("H134").Value = 6722.8
$ws.Range("I134").Value = 2801.5
$ws.Range("J134").Value = 13694
$ws.Range("K134").Value = 8404.5
$ws.Range("L134").Value = 41082
$ws.Range("M134").Value = -5869.5
$ws.Range("N134").Value = -46152
$ws.Range("H135").Value = 98395.2
$ws.Range("J135").Value = 98395.2
$ws.Range("L135").Value = 98395.2
$ws.Range("N135").Value = -108535.2
$ws.Range("H138").Value = 88874.14
$ws.Range("J138").Value = 88874.14
$ws.Range("L138").Value = 88874.14
$ws.Range("N138").Value = -99154.14

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 205019
$ws.Range("I4").Value = 256262.5
$ws.Range("J4").Value = 45
$ws.Range("K4").Value = 256262.5
$ws.Range("L4").Value = 45
$ws.Range("N4").Value = -269
$ws.Range("H62").Value = 3849.25
$ws.Range("I62").Value = 3849.25
$ws.Range("K62").Value = 3849.25
$ws.Range("M62").Value = -3225.25
$ws.Range("H65").Value = 3849.25
$ws.Range("I65").Value = 3849.25
$ws.Range("K65").Value = 19246.25
$ws.Range("M65").Value = -16126.25
$ws.Range("H132").Value = 1300826.1
$ws.Range("J132").Value = 2745.6667
$ws.Range("L132").Value = 8237.000100000001
$ws.Range("N132").Value = -13297.0001
$ws.Range("H138").Value = 51631.555
$ws.Range("J138").Value = 51631.555
$ws.Range("L138").Value = 51631.555
$ws.Range("N138").Value = -61911.555
$ws.Range("M4").Value = -256150.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 5240.136
$ws.Range("J132").Value = 6328.8823
$ws.Range("L132").Value = 56959.9407
$ws.Range("N132").Value = -62019.9407

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 14999.4
$ws.Range("J47").Value = 14999.4
$ws.Range("L47").Value = 14999.4
$ws.Range("N47").Value = -16135.4
$ws.Range("H93").Value = 18363
$ws.Range("J93").Value = 18363
$ws.Range("L93").Value = 18363
$ws.Range("N93").Value = -22107
$ws.Range("H135").Value = 69998.336
$ws.Range("J135").Value = 69998.336
$ws.Range("L135").Value = 69998.336
$ws.Range("N135").Value = -80138.336
$ws.Range("H140").Value = 87372.5
$ws.Range("J140").Value = 86925.71000000001
$ws.Range("L140").Value = 86925.71000000001
$ws.Range("N140").Value = -97285.71000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2660.7693
$ws.Range("I16").Value = 2621.818
$ws.Range("K16").Value = 2621.818
$ws.Range("M16").Value = -2451.818
$ws.Range("H46").Value = 2684.8215
$ws.Range("J46").Value = 3404.1667
$ws.Range("L46").Value = 3404.1667
$ws.Range("N46").Value = -3780.1667
$ws.Range("H132").Value = 3396
$ws.Range("I132").Value = 3154.4
$ws.Range("K132").Value = 9463.200000000001
$ws.Range("M132").Value = -6933.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value = 60390
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H132").Value = 1465.7
$ws.Range("I132").Value = 1305.1305
$ws.Range("K132").Value = 3915.3915
$ws.Range("M132").Value = -1385.3915
$ws.Range("N127").ClearContents()
